# fix(module3): use uncon_planned_qty for future production; keep produced for today
# Updates DeliveryPlan (ori_deployment_uid / material / delivery_qty / VFR) and
# VehicleLog (total_units / total_volume / VFR) to reflect the corrected
# future-production read (uncon_planned_qty) vs. today's produced qty.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet: DeliveryPlan
# Columns: A vehicle_uid | B ori_deployment_uid | C material | D sending
#          E receiving | F planned_deployment_date | G actual_ship_date
#          H actual_delivery_date | I delivery_qty | J truck_type
#          K truck_load_pct | L WFR | M VFR
# ---------------------------------------------------------------------
$dp = $wb.Worksheets.Item("DeliveryPlan")

# Row 2
$dp.Cells.Item(2, 2).Value = "MAT_A|PLANT_001|DC_001|2024-01-02|net demand for forecast|000014"
$dp.Cells.Item(2, 3).Value = "MAT_A"
$dp.Cells.Item(2, 9).Value = 20
$dp.Cells.Item(2, 13).Value = 0.9625

# Row 3
$dp.Cells.Item(3, 2).Value = "MAT_B|PLANT_001|DC_001|2024-01-02|net demand for forecast|000019"
$dp.Cells.Item(3, 3).Value = "MAT_B"
$dp.Cells.Item(3, 9).Value = 40
$dp.Cells.Item(3, 13).Value = 0.9625

# Row 4
$dp.Cells.Item(4, 2).Value = "MAT_A|PLANT_001|DC_001|2024-01-03|net demand for forecast|000015"
$dp.Cells.Item(4, 3).Value = "MAT_A"
$dp.Cells.Item(4, 9).Value = 20
$dp.Cells.Item(4, 13).Value = 0.9625

# Row 5
$dp.Cells.Item(5, 2).Value = "MAT_B|PLANT_001|DC_001|2024-01-03|net demand for forecast|000020"
$dp.Cells.Item(5, 3).Value = "MAT_B"
$dp.Cells.Item(5, 9).Value = 35
$dp.Cells.Item(5, 13).Value = 0.9625

# Row 6
$dp.Cells.Item(6, 2).Value = "MAT_B|PLANT_001|DC_001|2024-01-03|net demand for forecast|000020"
$dp.Cells.Item(6, 3).Value = "MAT_B"
$dp.Cells.Item(6, 9).Value = 5
$dp.Cells.Item(6, 13).Value = 0.9625

# Row 7
$dp.Cells.Item(7, 2).Value = "MAT_A|PLANT_001|DC_001|2024-01-04|net demand for forecast|000016"
$dp.Cells.Item(7, 3).Value = "MAT_A"
$dp.Cells.Item(7, 9).Value = 20
$dp.Cells.Item(7, 13).Value = 0.9625

# Row 8
$dp.Cells.Item(8, 2).Value = "MAT_B|PLANT_001|DC_001|2024-01-04|net demand for forecast|000021"
$dp.Cells.Item(8, 3).Value = "MAT_B"
$dp.Cells.Item(8, 9).Value = 40
$dp.Cells.Item(8, 13).Value = 0.9625

# Row 9
$dp.Cells.Item(9, 2).Value = "MAT_A|PLANT_001|DC_001|2024-01-05|net demand for forecast|000017"
$dp.Cells.Item(9, 3).Value = "MAT_A"
$dp.Cells.Item(9, 9).Value = 20
$dp.Cells.Item(9, 13).Value = 0.9625

# Row 10
$dp.Cells.Item(10, 2).Value = "MAT_B|PLANT_001|DC_001|2024-01-05|net demand for forecast|000022"
$dp.Cells.Item(10, 3).Value = "MAT_B"
$dp.Cells.Item(10, 9).Value = 30
$dp.Cells.Item(10, 13).Value = 0.9625

# Row 11 (material unchanged: MAT_A, only the underlying deployment uid changes)
$dp.Cells.Item(11, 2).Value = "MAT_A|PLANT_001|DC_002|2024-01-02|net demand for forecast|000010"
$dp.Cells.Item(11, 3).Value = "MAT_A"

# Row 12
$dp.Cells.Item(12, 2).Value = "MAT_A|PLANT_001|DC_002|2024-01-03|net demand for forecast|000011"
$dp.Cells.Item(12, 3).Value = "MAT_A"

# Row 13
$dp.Cells.Item(13, 2).Value = "MAT_A|PLANT_001|DC_002|2024-01-03|net demand for forecast|000011"
$dp.Cells.Item(13, 3).Value = "MAT_A"

# Row 14
$dp.Cells.Item(14, 2).Value = "MAT_A|PLANT_001|DC_002|2024-01-04|net demand for forecast|000012"
$dp.Cells.Item(14, 3).Value = "MAT_A"

# ---------------------------------------------------------------------
# Sheet: VehicleLog
# Columns: A date | B sending | C receiving | D truck_type | E vehicle_no
#          F vehicle_uid | G total_units | H total_weight | I total_volume
#          J WFR | K VFR | L truck_used
# ---------------------------------------------------------------------
$vl = $wb.Worksheets.Item("VehicleLog")

# Row 2 (vehicle #1)
$vl.Cells.Item(2, 7).Value = 115
$vl.Cells.Item(2, 9).Value = 192.5
$vl.Cells.Item(2, 11).Value = 0.9625

# Row 3 (vehicle #2)
$vl.Cells.Item(3, 7).Value = 115
$vl.Cells.Item(3, 9).Value = 192.5
$vl.Cells.Item(3, 11).Value = 0.9625
